# feat(doc): add space in templates
#
# Insert a new, empty paragraph right after the "Генеральный директор"
# paragraph (and before the "___________________ Михайлов Д. С." /
# signature paragraph) in the left-hand signature block of the table.
# The new paragraph carries the same run-level formatting (sz/szCs 20)
# as its neighbours and disables widow/orphan control, matching the
# surrounding paragraphs in this template, but contains no runs/text.

$d = $word.ActiveDocument

# Locate the "Генеральный директор" paragraph text via Find (as in the
# rest of the template, this string only occurs once, right above the
# signature line).
$rng = $d.Content
$found = $rng.Find.Execute("Генеральный директор", $true, $false, $false,
                            $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'Генеральный директор' paragraph"
}

# Collapse the found range to its end (right after the word
# "директор", i.e. still inside that paragraph, just before its
# paragraph mark) so the insert below only adds new content and leaves
# the existing paragraph completely untouched.
$rng.Collapse(0)

# A single empty paragraph, formatted like its neighbours:
# widowControl off, run-properties sz/szCs = 20 half-points (10pt), and
# no runs at all (a fully empty paragraph).
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr>' +
    '<w:widowControl w:val="0"/>' +
    '<w:rPr>' +
    '<w:sz w:val="20"/>' +
    '<w:szCs w:val="20"/>' +
    '</w:rPr>' +
    '</w:pPr>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# InsertXML on a collapsed (zero-length) range inserts the supplied
# OOXML right at that point without touching anything before/after it,
# which splits a brand-new paragraph in after "Генеральный директор".
$rng.InsertXML($newParaXml)
